$wb = $excel.ActiveWorkbook

# --- Sheet1: add the new "header4" column header (D1), bold like the other headers ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("D1").Value = "header4"
$ws1.Range("D1").Font.Bold = $true

# --- Make Sheet1 the active sheet/tab (previously "2ndsheet" was active) and
#     move its selection to D2. "2ndsheet" keeps its own previous selection (B3)
#     and simply stops being the active tab once Sheet1 is activated. ---
$ws1.Activate()
$ws1.Range("D2").Select()
